$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "prefix|Emissions|BC|Harmonized"
$ws.Range("D3").Value = "prefix|Emissions|Sulfur|Harmonized"
